# Auto-generated from the cryptos.xlsx diff: update Price (D) and
# Volume(1h) (E) columns for rows 2-51, plus the VeChain/FraxShare
# row swap (rows 44-45: Coin name + Link + Price + Volume).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.879.34"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.213.62"
$ws.Range("E3").Value = "  -0.21%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'291.92"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'87.15"
$ws.Range("E6").Value = "  +1.14%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.467"
$ws.Range("E9").Value = "  -0.86%  "
$ws.Range("D10").Value = "'30.42"
$ws.Range("E10").Value = "  -0.97%  "
$ws.Range("E11").Value = "  -0.69%  "
$ws.Range("D12").Value = "'49.99"
$ws.Range("E12").Value = "  +5.16%  "
$ws.Range("E13").Value = "  +2.42%  "
$ws.Range("D14").Value = "'6.44"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "2.554.83"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("D16").Value = "'13.76"
$ws.Range("E16").Value = "  -1.95%  "
$ws.Range("D17").Value = "2.204.92"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "39.812.92"
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").Value = "'11.13"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").Value = "'5.73"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'237.00"
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'2.45"
$ws.Range("E26").Value = "  -0.53%  "
$ws.Range("D27").Value = "'1.82"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "'23.13"
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("D30").Value = "'2.04"
$ws.Range("E30").Value = "  -7.13%  "
$ws.Range("D31").Value = "'157.11"
$ws.Range("E31").Value = "  +3.63%  "
$ws.Range("D32").Value = "'31.89"
$ws.Range("E32").Value = "  -3.25%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").Value = "'0.0709"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("E37").Value = "  -1.60%  "
$ws.Range("D39").Value = "'0.0983"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'1.73"
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("D41").Value = "'15.24"
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("D42").Value = "2.112.43"
$ws.Range("E42").Value = "  +2.19%  "
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0269"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.99"
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").Value = "'17.74"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'2.08"
$ws.Range("E47").Value = "  -0.58%  "
$ws.Range("D48").Value = "'2.68"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "2.425.34"
$ws.Range("E49").Value = "  -0.33%  "
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("D51").Value = "'88.50"
$ws.Range("E51").Value = "  -0.57%  "
